# kundur_reg.xlsx :: REPCA1 sheet
# `Model.s_update` now calls the lambdified and numerical function for each
# service individually, so REPCA1 needs three new per-service flag columns:
# VCFlag, RefFlag, Fflag. Insert them right after the existing "busf" column
# (column H) and before "Tfltr" (old column I), shifting the rest of the
# sheet's columns to the right, then seed the new columns' single data row
# with 0 (matching the other flag-style columns on this sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("REPCA1")

# Insert 3 new blank columns at I:K, pushing existing I:AH to L:AK.
$ws.Columns("I:K").Insert()

# Header row (row 1).
$ws.Range("I1").Value = "VCFlag"
$ws.Range("J1").Value = "RefFlag"
$ws.Range("K1").Value = "Fflag"

# Single data row (row 2) for the REPCA1_1 device.
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
